$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 5 new rows (8-12) below the existing data row 7. This pushes the
# old totals row (old row 8) and footer row (old row 9) down to rows 13/14,
# matching the diff's row shift.
$ws.Rows("8:12").Insert()

# Clone row 7's full formatting (styles, merges) into each new row so the
# new rows look identical to the template row.
$ws.Range("A7:Q7").Copy($ws.Range("A8:Q8"))
$ws.Range("A7:Q7").Copy($ws.Range("A9:Q9"))
$ws.Range("A7:Q7").Copy($ws.Range("A10:Q10"))
$ws.Range("A7:Q7").Copy($ws.Range("A11:Q11"))
$ws.Range("A7:Q7").Copy($ws.Range("A12:Q12"))

# Restore the per-row heights from the target layout.
$ws.Rows("7").RowHeight = 25.5
$ws.Rows("8").RowHeight = 24.75
$ws.Rows("9").RowHeight = 25.5
$ws.Rows("10").RowHeight = 24.75
$ws.Rows("11").RowHeight = 25.5
$ws.Rows("12").RowHeight = 25.5

function Set-ItemRow($rowNum, $idx, $name, $balance, $orderLimit, $price, $sellPrice, $txCount) {
    $ws.Range("A$rowNum").Value = $idx

    $ws.Range("C$rowNum").NumberFormat = "@"
    $ws.Range("C$rowNum").Value = $name

    $ws.Range("H$rowNum").NumberFormat = "@"
    $ws.Range("H$rowNum").Value = $balance

    $ws.Range("L$rowNum").NumberFormat = "@"
    $ws.Range("L$rowNum").Value = $orderLimit

    $ws.Range("N$rowNum").NumberFormat = "@"
    $ws.Range("N$rowNum").Value = $price

    $ws.Range("P$rowNum").NumberFormat = "@"
    $ws.Range("P$rowNum").Value = $sellPrice

    $ws.Range("Q$rowNum").NumberFormat = "@"
    $ws.Range("Q$rowNum").Value = $txCount
}

Set-ItemRow 7  1 "BRIMOSALM EYE DROPS 5 ML"            "0:0"  "1" "41.00" "41.0000" "1:0"
Set-ItemRow 8  2 "CATAFAST 50 MG 9 SACHET"              "1:7"  "1" "72.00" "7.9200"  "0:1"
Set-ItemRow 9  3 "DECLOPHEN 75MG/3ML 3 AMPOULES"        "7:1"  "1" "36.00" "36.0000" "1:0"
Set-ItemRow 10 4 "DELTAVIT B12 1MG 30 SUBLINGUAL TAB"   "0:0"  "1" "75.00" "75.0000" "1:0"
Set-ItemRow 11 5 "NOSTAMINE EYE/NOSE DROPS 15 ML"       "0:0"  "1" "22.00" "22.0000" "1:0"
Set-ItemRow 12 6 "قطن 100 جم"                            "30:0" "0" "20.00" "20.0000" "1:0"

# Totals row (was row 8, now row 13): sum of the "sell price" column.
$ws.Range("P13").Value = 201.91999999999999

# Footer row (was row 9, now row 14): refreshed generation timestamp.
$ws.Range("A14").Value = "Thursday, 4 September, 2025 11:33 AM"
